$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("M21").Value = ""
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("M23").Value = ""
$ws.Range("H29").Value = 1934
$ws.Range("I29").Value = 3
$ws.Range("J29").Value = 2899.5
$ws.Range("K29").Value = 9
$ws.Range("L29").Value = 8698.5
$ws.Range("M29").Value = 272
$ws.Range("N29").Value = -9260.5
$ws.Range("H38").Value = 1021.6
$ws.Range("I38").Value = 1002
$ws.Range("J38").Value = 1100
$ws.Range("K38").Value = 3006
$ws.Range("L38").Value = 3300
$ws.Range("M38").Value = -2634
$ws.Range("N38").Value = -4044
$ws.Range("H58").Value = 1613.9
$ws.Range("J58").Value = 3050
$ws.Range("L58").Value = 9150
$ws.Range("N58").Value = -9450
$ws.Range("H62").Value = 7333
$ws.Range("I62").Value = 5999.5
$ws.Range("J62").Value = 10000
$ws.Range("K62").Value = 5999.5
$ws.Range("L62").Value = 10000
$ws.Range("M62").Value = -5375.5
$ws.Range("N62").Value = -11248
$ws.Range("H65").Value = 7333
$ws.Range("I65").Value = 5999.5
$ws.Range("J65").Value = 10000
$ws.Range("K65").Value = 29997.5
$ws.Range("L65").Value = 50000
$ws.Range("M65").Value = -26877.5
$ws.Range("N65").Value = -56240
$ws.Range("H112").Value = 2204.5386
$ws.Range("J112").Value = 2204.5386
$ws.Range("L112").Value = 6613.6158
$ws.Range("N112").Value = -8829.6158
$ws.Range("H132").Value = 2551.6924
$ws.Range("I132").Value = 2551.6924
$ws.Range("K132").Value = 7655.0772
$ws.Range("M132").Value = -5125.0772
$ws.Range("H137").Value = 2536.25
$ws.Range("I137").Value = 1893.9166
$ws.Range("J137").Value = 3499.75
$ws.Range("K137").Value = 5681.7498
$ws.Range("L137").Value = 10499.25
$ws.Range("M137").Value = -3131.7498
$ws.Range("N137").Value = -15599.25
$ws.Range("H138").Value = 4638.3335
$ws.Range("I138").Value = 1104.3636
$ws.Range("J138").Value = 7067.9375
$ws.Range("K138").Value = 3313.0908
$ws.Range("L138").Value = 21203.8125
$ws.Range("M138").Value = 1826.9092
$ws.Range("N138").Value = -31483.8125

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1838.2
$ws.Range("I45").Value = 1838.2
$ws.Range("K45").Value = 1838.2
$ws.Range("M45").Value = -1461.2
$ws.Range("H46").Value = 15190.25
$ws.Range("I46").Value = 16143.875
$ws.Range("J46").Value = 13283
$ws.Range("K46").Value = 16143.875
$ws.Range("L46").Value = 13283
$ws.Range("M46").Value = -15824.875
$ws.Range("N46").Value = -13921
$ws.Range("H52").Value = 29999
$ws.Range("J52").Value = 29999
$ws.Range("L52").Value = 29999
$ws.Range("N52").Value = -30635
$ws.Range("H61").Value = 1511.4688
$ws.Range("I61").Value = 1345.4286
$ws.Range("J61").Value = 2673.75
$ws.Range("K61").Value = 1345.4286
$ws.Range("L61").Value = 2673.75
$ws.Range("M61").Value = -1133.4286
$ws.Range("N61").Value = -3097.75
$ws.Range("H110").Value = 7744.727
$ws.Range("I110").Value = 5574.125
$ws.Range("K110").Value = 5574.125
$ws.Range("M110").Value = -3529.125
$ws.Range("H132").Value = 1684
$ws.Range("I132").Value = 1584.1177
$ws.Range("J132").Value = 2250
$ws.Range("K132").Value = 4752.3531
$ws.Range("L132").Value = 6750
$ws.Range("M132").Value = -2222.3531
$ws.Range("N132").Value = -11810
$ws.Range("H136").Value = 1511.4688
$ws.Range("I136").Value = 1345.4286
$ws.Range("J136").Value = 2673.75
$ws.Range("K136").Value = 4036.2858
$ws.Range("L136").Value = 8021.25
$ws.Range("M136").Value = -1486.2858
$ws.Range("N136").Value = -13121.25
$ws.Range("H140").Value = 84747.5
$ws.Range("J140").Value = 91330
$ws.Range("L140").Value = 91330
$ws.Range("N140").Value = -101690

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2717.5789
$ws.Range("I105").Value = 2590.7778
$ws.Range("K105").Value = 2590.7778
$ws.Range("M105").Value = -843.7777999999998
$ws.Range("H107").Value = 1237.2
$ws.Range("I107").Value = 1246.5
$ws.Range("K107").Value = 1246.5
$ws.Range("M107").Value = 673.5
$ws.Range("H134").Value = 2828.5173
$ws.Range("I134").Value = 2685
$ws.Range("J134").Value = 3147.4443
$ws.Range("K134").Value = 8055
$ws.Range("L134").Value = 9442.332900000001
$ws.Range("M134").Value = -5520
$ws.Range("N134").Value = -14512.3329

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4634.4316
$ws.Range("I31").Value = 2826.6
$ws.Range("K31").Value = 2826.6
$ws.Range("M31").Value = -2531.6
$ws.Range("H34").Value = 4634.4316
$ws.Range("I34").Value = 2826.6
$ws.Range("K34").Value = 2826.6
$ws.Range("M34").Value = -2624.6
$ws.Range("H54").Value = 19197.334
$ws.Range("J54").Value = 19197.334
$ws.Range("L54").Value = 19197.334
$ws.Range("N54").Value = -20513.334
$ws.Range("H58").Value = 3042.4827
$ws.Range("I58").Value = 1348.9375
$ws.Range("K58").Value = 1348.9375
$ws.Range("M58").Value = -1145.9375
$ws.Range("H86").Value = 10198.4
$ws.Range("I86").Value = 8864.166999999999
$ws.Range("K86").Value = 8864.166999999999
$ws.Range("M86").Value = -7741.166999999999
$ws.Range("H89").Value = 10198.4
$ws.Range("I89").Value = 8864.166999999999
$ws.Range("K89").Value = 44320.835
$ws.Range("M89").Value = -38704.835
$ws.Range("H105").Value = 1242.3636
$ws.Range("I105").Value = 1030.875
$ws.Range("J105").Value = 1806.3334
$ws.Range("K105").Value = 1030.875
$ws.Range("L105").Value = 1806.3334
$ws.Range("M105").Value = 716.125
$ws.Range("N105").Value = -5300.3334
$ws.Range("H132").Value = 2441.372
$ws.Range("I132").Value = 1843.1892
$ws.Range("J132").Value = 6130.1665
$ws.Range("K132").Value = 5529.5676
$ws.Range("L132").Value = 18390.4995
$ws.Range("M132").Value = -2999.5676
$ws.Range("N132").Value = -23450.4995
$ws.Range("H134").Value = 2278.8723
$ws.Range("I134").Value = 1813.1875
$ws.Range("J134").Value = 3272.3333
$ws.Range("K134").Value = 5439.5625
$ws.Range("L134").Value = 9816.999899999999
$ws.Range("M134").Value = -2904.5625
$ws.Range("N134").Value = -14886.9999
$ws.Range("H136").Value = 3042.4827
$ws.Range("I136").Value = 1348.9375
$ws.Range("K136").Value = 4046.8125
$ws.Range("M136").Value = -1496.8125

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 25070210
$ws.Range("I4").Value = 30635782
$ws.Range("J4").Value = 25135.375
$ws.Range("K4").Value = 91907346
$ws.Range("L4").Value = 75406.125
$ws.Range("M4").Value = -91907234
$ws.Range("N4").Value = -75630.125
$ws.Range("H17").Value = 1374.3334
$ws.Range("J17").Value = 2698.6667
$ws.Range("L17").Value = 8096.000100000001
$ws.Range("N17").Value = -8434.000100000001
$ws.Range("H92").Value = 394.68182
$ws.Range("I92").Value = 401.5
$ws.Range("J92").Value = 389
$ws.Range("K92").Value = 1204.5
$ws.Range("L92").Value = 1167
$ws.Range("M92").Value = 43.5
$ws.Range("N92").Value = -3663
$ws.Range("H107").Value = 402.16666
$ws.Range("I107").Value = 278.85715
$ws.Range("K107").Value = 836.5714499999999
$ws.Range("M107").Value = 1083.42855

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6653.353
$ws.Range("I70").Value = 6205.7
$ws.Range("K70").Value = 6205.7
$ws.Range("M70").Value = -5935.7
$ws.Range("H73").Value = 6653.353
$ws.Range("I73").Value = 6205.7
$ws.Range("K73").Value = 6205.7
$ws.Range("M73").Value = -5269.7
$ws.Range("H132").Value = 2361.9312
$ws.Range("I132").Value = 1816.5264
$ws.Range("K132").Value = 5449.5792
$ws.Range("M132").Value = -2919.5792

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H43").Value = 1553875
$ws.Range("J43").Value = 2373555.5
$ws.Range("L43").Value = 2373555.5
$ws.Range("N43").Value = -2373941.5
$ws.Range("H100").Value = 1405.5
$ws.Range("I100").Value = 933
$ws.Range("J100").Value = 1689
$ws.Range("K100").Value = 933
$ws.Range("L100").Value = 1689
$ws.Range("M100").Value = -392
$ws.Range("N100").Value = -2771
$ws.Range("H121").Value = 0
$ws.Range("J121").Value = 0
$ws.Range("L121").Value = 0
$ws.Range("N121").Value = ""
$ws.Range("H132").Value = 3390.375
$ws.Range("I132").Value = 2986.6956
$ws.Range("K132").Value = 8960.086800000001
$ws.Range("M132").Value = -6430.086800000001
$ws.Range("H136").Value = 5077.2666
$ws.Range("I136").Value = 4855.5835
$ws.Range("K136").Value = 14566.7505
$ws.Range("M136").Value = -12016.7505

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 500
$ws.Range("I2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("M2").Value = ""
$ws.Range("H29").Value = 26350
$ws.Range("J29").Value = 7700
$ws.Range("L29").Value = 7700
$ws.Range("N29").Value = -8280
$ws.Range("H107").Value = 749.3333
$ws.Range("I107").Value = 874
$ws.Range("K107").Value = 2622
$ws.Range("M107").Value = -702
$ws.Range("H122").Value = 3069.5
$ws.Range("I122").Value = 3069.5
$ws.Range("K122").Value = 9208.5
$ws.Range("M122").Value = -6758.5
$ws.Range("H132").Value = 60104.438
$ws.Range("I132").Value = 73528.84
$ws.Range("K132").Value = 220586.52
$ws.Range("M132").Value = -218056.52
